# Weekly price-sheet update: a new week's Mango price record is inserted
# right after the header/first block, at row 40. Every existing data row
# from 40 down to 150 shifts down by one (to 41..151); the sheet's used
# range grows from A1:T150 to A1:T151.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 40 (pushes rows 40:150 -> 41:151).
$ws.Rows("40:40").Insert()

# Populate the newly inserted row 40 with this week's record.
$fecha = Get-Date -Year 2023 -Month 2 -Day 8 -Hour 0 -Minute 0 -Second 0

$ws.Range("A40").Value = 11
$ws.Range("B40").Value = "Vega Monumental Concepción"
$ws.Range("C40").Value = "Bíobío"
$ws.Range("D40").Value = $fecha.Date
$ws.Range("E40").Value = 8
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100108
$ws.Range("H40").Value = "Tropicales y subtropicales"
$ws.Range("I40").Value = 100108002
$ws.Range("J40").Value = "Mango"
$ws.Range("K40").Value = "Sin especificar"
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 200
$ws.Range("N40").Value = 6500
$ws.Range("O40").Value = 7000
$ws.Range("P40").Value = 6750
$ws.Range("Q40").Value = "$/bandeja 4 kilos"
$ws.Range("R40").Value = "Perú"
$ws.Range("S40").Value = 1688
$ws.Range("T40").Value = 4
